$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for new columns I and J, matching the style of existing headers (H1 etc.)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Values for column I (I0) and column J (IF), rows 2-18
$values = @(
    @{ Row = 2;  I = 8; J = 9 },
    @{ Row = 3;  I = 7; J = 7 },
    @{ Row = 4;  I = 8; J = 8 },
    @{ Row = 5;  I = 6; J = 6 },
    @{ Row = 6;  I = 8; J = 8 },
    @{ Row = 7;  I = 7; J = 7 },
    @{ Row = 8;  I = 1; J = 4 },
    @{ Row = 9;  I = 1; J = 4 },
    @{ Row = 10; I = 1; J = 2 },
    @{ Row = 11; I = 1; J = 6 },
    @{ Row = 12; I = 1; J = 5 },
    @{ Row = 13; I = 1; J = 7 },
    @{ Row = 14; I = 1; J = 6 },
    @{ Row = 15; I = 1; J = 6 },
    @{ Row = 16; I = 1; J = 5 },
    @{ Row = 17; I = 1; J = 4 },
    @{ Row = 18; I = 4; J = 6 }
)

foreach ($entry in $values) {
    $ws.Cells.Item($entry.Row, 9).Value = $entry.I
    $ws.Cells.Item($entry.Row, 10).Value = $entry.J
}
